$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 6.449
$ws.Range("E4").Value = 12.682
$ws.Range("E5").Value = 13.254
$ws.Range("B7").Value = 6.956999999999999
$ws.Range("E8").Value = 13.718
$ws.Range("B16").Value = 6.782999999999999
$ws.Range("E16").Value = 12.914
